$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Codigo Operadora" in BI1, right after the existing
# "Data" header in BH1.
$ws.Range("BI1").Value = "Codigo Operadora"
